$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cell updates mirroring the refreshed crypto price/volume scrape.
# D (Price) and E (Volume(1h)) columns hold numeric-looking text, so the
# cell is forced to Text format before the value is written - this keeps
# Excel from auto-converting "255.66" / "0.17%" into a real number/percent.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "255.66"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.17%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.98"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-3.97%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.612"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-11.33%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05887"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.40%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.639"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.97%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8680"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9451"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.49%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1400"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.80%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.03778"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "9.10%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07069"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.27%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03206"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.22%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09248"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.43%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001537"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.46%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006020"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.91%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006015"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.96%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.516"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.46%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.193"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.59%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-3.25%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.04%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "8.75%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04216"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.06%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.26%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001222"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.05%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004270"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-6.59%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.07%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001507"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "2.90%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03808"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.28%"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1097"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.47%"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003916"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-30.57%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002277"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.46%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01155"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "18.63%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005464"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.48%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.02%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.06019"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-33.05%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002279"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "7.03%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.02%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.02%"
